# Auto-generated edit script: applies the targeted cell updates for rows 27-39
# of the "Artfynd" worksheet, matching the canonical OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27
$ws.Range("A27").Value = 111612720
$ws.Range("B27").Value = 88924
$ws.Range("D27").Value = "LC"
$ws.Range("E27").Value = 256703
$ws.Range("F27").Value = "Tallfingersvamp"
$ws.Range("G27").Value = "Ramaria eosanguinea"
$ws.Range("H27").Value = "R.H.Petersen"
$ws.Range("I27").Value = "'2"
$ws.Range("J27").Value = ""
$ws.Range("L27").Value = ""
$ws.Range("M27").Value = ""
$ws.Range("Q27").Value = 491993.9996831641
$ws.Range("R27").Value = 6785505.377163783
$ws.Range("S27").Value = 100
$ws.Range("AF27").Value = ""

# Row 28
$ws.Range("A28").Value = 111612738
$ws.Range("B28").Value = 56414
$ws.Range("E28").Value = 100049
$ws.Range("F28").Value = "Spillkråka"
$ws.Range("G28").Value = "Dryocopus martius"

# Row 29
$ws.Range("A29").Value = 111612726
$ws.Range("B29").Value = 90168
$ws.Range("D29").Value = "VU"
$ws.Range("E29").Value = 717
$ws.Range("F29").Value = "Borsttagging"
$ws.Range("G29").Value = "Gloiodon strigosus"
$ws.Range("H29").Value = "(Schwein. : Fr.) P. Karst."
$ws.Range("I29").Value = ""
$ws.Range("Q29").Value = 491952.3910193561
$ws.Range("R29").Value = 6785464.984647369
$ws.Range("S29").Value = 10

# Row 30
$ws.Range("A30").Value = 111612736
$ws.Range("B30").Value = 56398
$ws.Range("D30").Value = "NT"
$ws.Range("E30").Value = 100109
$ws.Range("F30").Value = "Tretåig hackspett"
$ws.Range("G30").Value = "Picoides tridactylus"
$ws.Range("H30").Value = "(Linnaeus, 1758)"
$ws.Range("J30").Value = ""
$ws.Range("L30").Value = ""
$ws.Range("M30").Value = "äldre spår"
$ws.Range("AF30").Value = ""

# Row 31
$ws.Range("A31").Value = 111682655
$ws.Range("B31").Value = 90658
$ws.Range("E31").Value = 4361
$ws.Range("F31").Value = "Orange taggsvamp"
$ws.Range("G31").Value = "Hydnellum aurantiacum"
$ws.Range("H31").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("I31").Value = "'20"

# Row 32
$ws.Range("A32").Value = 111682652
$ws.Range("B32").Value = 90709
$ws.Range("E32").Value = 5448
$ws.Range("F32").Value = "Svartvit taggsvamp"
$ws.Range("G32").Value = "Phellodon connatus"
$ws.Range("H32").Value = "(Schultz) nom.prov"
$ws.Range("I32").Value = ""

# Row 33
$ws.Range("A33").Value = 111682658
$ws.Range("B33").Value = 90689
$ws.Range("E33").Value = 5966
$ws.Range("F33").Value = "Motaggsvamp"
$ws.Range("G33").Value = "Sarcodon squamosus"
$ws.Range("H33").Value = "(Schaeff.) Quél."

# Row 34
$ws.Range("A34").Value = 111682665
$ws.Range("B34").Value = 90682
$ws.Range("E34").Value = 2059
$ws.Range("F34").Value = "Skrovlig taggsvamp"
$ws.Range("G34").Value = "Hydnellum scabrosum"
$ws.Range("H34").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"

# Row 36
$ws.Range("A36").Value = 112045343
$ws.Range("B36").Value = 90658
$ws.Range("D36").Value = "NT"
$ws.Range("E36").Value = 4361
$ws.Range("F36").Value = "Orange taggsvamp"
$ws.Range("G36").Value = "Hydnellum aurantiacum"
$ws.Range("H36").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("I36").Value = "'10"
$ws.Range("N36").Value = ""
$ws.Range("P36").Value = "Nedre Tetvasseltjärnen (Nedre Tetvasseltjärnen), Dlr"
$ws.Range("Q36").Value = 492044.3381435904
$ws.Range("R36").Value = 6785564.065369682
$ws.Range("Z36").Value = "12:30"
$ws.Range("AB36").Value = "12:30"
$ws.Range("AF36").Value = ""
$ws.Range("AH36").Value = ""
$ws.Range("AJ36").Value = ""
$ws.Range("AK36").Value = ""
$ws.Range("AO36").Value = ""
$ws.Range("AW36").Value = "Bo karlstens"
$ws.Range("AX36").Value = "Bo karlstens, Bengt Oldhammer, Janolof Hermansson, Birgitta Kvist"

# Row 37
$ws.Range("A37").Value = 112073422
$ws.Range("B37").Value = 88949
$ws.Range("D37").Value = "LC"
$ws.Range("E37").Value = 233195
$ws.Range("F37").Value = ""
$ws.Range("G37").Value = "Ramaria neoformosa"
$ws.Range("H37").Value = "sensu Schild"
$ws.Range("I37").Value = "'3"
$ws.Range("J37").Value = "fruktkroppar"
$ws.Range("N37").Value = ""
$ws.Range("P37").Value = "N om Nedre Tetvasseltjärnen, Dlr"
$ws.Range("Q37").Value = 491995.9899496675
$ws.Range("R37").Value = 6785531.400109125
$ws.Range("S37").Value = 5
$ws.Range("Z37").Value = "00:00"
$ws.Range("AB37").Value = "00:00"
$ws.Range("AF37").Value = ""
$ws.Range("AH37").Value = "Sandtallskog"
$ws.Range("AJ37").Value = "tall"
$ws.Range("AK37").Value = "Pinus sylvestris"
$ws.Range("AO37").Value = "Pinus sylvestris"
$ws.Range("AW37").Value = "Janolof Hermansson"
$ws.Range("AX37").Value = "Janolof Hermansson, Bengt Oldhammer, Bo karlstens, Birgitta Kvist"

# Row 38
$ws.Range("A38").Value = 112045406
$ws.Range("B38").Value = 90682
$ws.Range("D38").Value = "NT"
$ws.Range("E38").Value = 2059
$ws.Range("F38").Value = "Skrovlig taggsvamp"
$ws.Range("G38").Value = "Hydnellum scabrosum"
$ws.Range("H38").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("I38").Value = ""
$ws.Range("J38").Value = ""
$ws.Range("N38").Value = ""
$ws.Range("P38").Value = "Nedre Tetvasseltjärnen (Nedre Tetvasseltjärnen), Dlr"
$ws.Range("Q38").Value = 492044.3381435904
$ws.Range("R38").Value = 6785564.065369682
$ws.Range("S38").Value = 15
$ws.Range("Z38").Value = "12:36"
$ws.Range("AB38").Value = "12:36"
$ws.Range("AF38").Value = ""
$ws.Range("AH38").Value = ""
$ws.Range("AJ38").Value = ""
$ws.Range("AK38").Value = ""
$ws.Range("AO38").Value = ""
$ws.Range("AW38").Value = "Bo karlstens"
$ws.Range("AX38").Value = "Bo karlstens, Bengt Oldhammer, Janolof Hermansson, Birgitta Kvist"

# Row 39
$ws.Range("A39").Value = 112073630
$ws.Range("B39").Value = 88924
$ws.Range("D39").Value = "LC"
$ws.Range("E39").Value = 256703
$ws.Range("F39").Value = "Tallfingersvamp"
$ws.Range("G39").Value = "Ramaria eosanguinea"
$ws.Range("H39").Value = "R.H.Petersen"
$ws.Range("I39").Value = "'1"
$ws.Range("J39").Value = "fruktkroppar"
$ws.Range("N39").Value = ""
$ws.Range("P39").Value = "N om Nedre Tetvasseltjärnen, Dlr"
$ws.Range("Q39").Value = 491917.2246397196
$ws.Range("R39").Value = 6785497.359069696
$ws.Range("Z39").Value = "00:00"
$ws.Range("AB39").Value = "00:00"
$ws.Range("AF39").Value = ""
$ws.Range("AH39").Value = "Sandtallskog"
$ws.Range("AJ39").Value = "tall"
$ws.Range("AK39").Value = "Pinus sylvestris"
$ws.Range("AO39").Value = "Pinus sylvestris"
$ws.Range("AW39").Value = "Janolof Hermansson"
$ws.Range("AX39").Value = "Janolof Hermansson, Bengt Oldhammer, Bo karlstens, Birgitta Kvist"
